# Update the "Reference ID" value in D2 (Furniture sheet) to a new
# numeric-looking reference id. The id must stay stored as text (it is a
# shared-string identifier, not a real number), so we build it in a
# scratch cell using a leading apostrophe (forces text type), copy just
# the *value* over with PasteSpecial (so no explicit cell style/number
# format is attached to D2 itself, matching the original cell which also
# carries no style), then clean up the scratch cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Value = "'41655678"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
